# Adds a new weekly price record for "Zanahoria" (Vega Monumental Concepción)
# at row 208, pushing the existing rows 208:237 down to 209:238.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 208 (shifts 208:237 -> 209:238)
$ws.Rows.Item(208).Insert()

# Populate the new row 208 with the new record
$ws.Range("A208").Value = 11
$ws.Range("B208").Value = "Vega Monumental Concepción"
$ws.Range("C208").Value = "Bíobío"
$ws.Range("D208").Value = 44748
$ws.Range("E208").Value = 8
$ws.Range("F208").Value = 100114013
$ws.Range("G208").Value = "Zanahoria"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 250
$ws.Range("K208").Value = 8000
$ws.Range("L208").Value = 8500
$ws.Range("M208").Value = 8200
$ws.Range("N208").Value = "$/saco 20 kilos"
$ws.Range("O208").Value = "Región de Coquimbo"
$ws.Range("P208").Value = 410
$ws.Range("Q208").Value = 20
$ws.Range("R208").Value = "Hortaliza"

# Match the date-number-format style used by the other D-column cells
$ws.Range("D208").Style = $ws.Range("D209").Style
